$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 922
$ws.Range("A3").Value = 1666
$ws.Range("A4").Value = 1667
$ws.Range("A5").Value = 4346
$ws.Range("A6").Value = 6932
$ws.Range("A7").Value = 6933
$ws.Range("A8").Value = 9701
